# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.711.51'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '2.370.06'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -4.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.88'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.09%  '
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.01'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -6.26%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '2.729.05'
$ws.Range("E14").Value = '  +2.98%  '
$ws.Range("D15").Value = '2.392.27'
$ws.Range("E15").Value = '  +3.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.814'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.64'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").Value = '45.622.42'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.06%  '
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.65'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.88'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.75'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.01%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  -12.82%  '
$ws.Range("E28").Value = '  -3.79%  '
$ws.Range("E29").Value = '  -3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.77'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +16.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.84'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.75'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.36'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.46'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0762'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.55%  '
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("E37").Value = '  +5.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.116'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.11'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.05%  '
$ws.Range("E40").Value = '  -7.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0294'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.19%  '
$ws.Range("D43").Value = '1.938.73'
$ws.Range("E43").Value = '  +4.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.69'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.25%  '
$ws.Range("E46").Value = '  -10.76%  '
$ws.Range("E47").Value = '  +7.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.30'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("E49").Value = '  -7.07%  '
$ws.Range("D50").Value = '2.599.27'
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '68.40'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.60%  '
